$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Recruitment Officer"
$ws.Range("A3").Value = "Personnel Officer"
$ws.Range("A4").Value = "General Affair Officer"
$ws.Range("A5").Value = "Trainner Officer"
$ws.Range("A6").Value = "Human Resource Supervisor"
$ws.Range("A7").Value = "Legal Supervisor"

$ws.Range("F8:F9").Select()
